# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K") previously stored a different, now-retired stat
# ("Strike#"). It has been regenerated from the per-game source data so
# that it now holds the true K (strikeout) counts. The freshly
# calculated s_vals below (one per game row, sheet rows 2-81) are
# written back into column G, replacing the stale values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly calculated K values (s_vals), in row order (row 2 .. row 81)
$s_vals = @(
    1,3,0,3,2,0,2,0,1,1,
    3,1,2,1,1,1,0,1,1,0,
    1,3,2,2,2,0,0,0,0,2,
    4,1,1,2,3,2,2,1,2,1,
    2,1,2,1,2,2,1,0,2,3,
    3,1,2,3,0,0,2,1,1,5,
    1,2,2,0,2,2,1,2,1,1,
    1,1,0,1,0,3,2,0,1,1
)

$firstRow = 2
$col = 7  # column G

for ($i = 0; $i -lt $s_vals.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, $col).Value = $s_vals[$i]
}
